$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document:
#       <empty run><b>Meta description</b>: Read our in-depth review...
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start
$boldLabel = "Meta description"
$restText = ": Read our in-depth review of Ali Baba's Gold and play for free. Enjoy stunning graphics, special features, and flexible interface for seamless gameplay."

$insertPoint = $d.Range($metaStart, $metaStart)
$insertPoint.InsertAfter($boldLabel + $restText)

# Bold only the "Meta description" label, leave the rest as-is.
$boldRange = $d.Range($metaStart, $metaStart + $boldLabel.Length)
$boldRange.Font.Bold = 1

# ------------------------------------------------------------------
# 2) At the end of the document: delete the bold "Play Ali Baba's
#    Gold Free - Review of Top Slot Machines" paragraph entirely
#    (it was duplicated there), and replace the text of the
#    following italic paragraph with the new feature-image prompt
#    (its italic formatting / existing leading empty run is left
#    untouched since we only rewrite the text content).
# ------------------------------------------------------------------
$oldBoldText = "Play Ali Baba's Gold Free - Review of Top Slot Machines"
$oldItalicText = "Read our in-depth review of Ali Baba's Gold and play for free. Enjoy stunning graphics, special features, and flexible interface for seamless gameplay."

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$secondLastPara = $d.Paragraphs.Item($count - 1)

if ($secondLastPara.Range.Text.Trim() -eq $oldBoldText) {
    $secondLastPara.Range.Delete()
}

$count = $d.Paragraphs.Count
$italicPara = $d.Paragraphs.Item($count)

$newImageText = "Create a feature image for ""Ali Baba's Gold"" that showcases a happy Maya warrior with glasses. The image should be in a cartoon style that captures the adventurous spirit of the game. The Maya warrior should be holding a bag of gold with a big smile on their face, indicating a successful treasure hunt. In the background, we can see the mysterious cave and the beautiful princess, along with the ruthless Jafar. The colors used should be bright and vibrant, perfectly capturing the essence of this beautiful online slot game. Overall, the image should be fun and engaging, inviting players to join Ali Baba on his quest for gold and love."

$replaceRange = $d.Range($italicPara.Range.Start, $italicPara.Range.End - 1)
$replaceRange.Text = $newImageText
